$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("I11").Value = "July"

# Row 12
$ws.Range("A12").Value = "STL526_539_537 | Sintra: Moorish Castle & Quinta da Regaleira e-Tickets"
$ws.Range("B12").Value = 8
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 8
$ws.Range("G12").Value = 24
$ws.Range("I12").Value = "June"

# Row 13
$ws.Range("A13").Value = "STL526_539_537 | Quinta da Regaleira & Moorish Castle Tickets with 3 Audios"
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = 8
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 12
$ws.Range("I13").Value = "August"

# Row 21
$ws.Range("A21").Value = "STL526_539_537 | Quinta da Regaleira & Moorish Castle Tickets with 3 Audios"
$ws.Range("B21").Value = 4
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 12
$ws.Range("I21").Value = "June"

# Row 22
$ws.Range("A22").Value = "STL526_539_537 | Sintra: Moorish Castle & Quinta da Regaleira e-Tickets"
$ws.Range("B22").Value = 8
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = 8
$ws.Range("G22").Value = 24
$ws.Range("I22").Value = "August"

# Row 26
$ws.Range("A26").Value = "STL539_537 | Quinta da Regaleira Ticket with Audio Guide & Sintra Tour"
$ws.Range("C26").Value = 0
$ws.Range("G26").Value = 6
$ws.Range("I26").Value = "May"

# Row 27
$ws.Range("A27").Value = "STL539_537 | Sintra: Quinta da Regaleira Ticket & Sintra Smartphone Tour"
$ws.Range("C27").Value = 6
$ws.Range("G27").Value = 12
$ws.Range("I27").Value = "August"

# Row 28
$ws.Range("A28").Value = "STL539_537 | Sintra: Quinta da Regaleira Ticket & Sintra Smartphone Tour"
$ws.Range("B28").Value = 6
$ws.Range("C28").Value = 6
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 12
$ws.Range("I28").Value = "June"

# Row 29
$ws.Range("A29").Value = "STL539_537 | Sintra: Quinta da Regaleira e-Ticket & Sintra City Audios"
$ws.Range("B29").Value = 0
$ws.Range("C29").Value = 0
$ws.Range("F29").Value = 6
$ws.Range("G29").Value = 6
$ws.Range("I29").Value = "July"

# Row 33
$ws.Range("I33").Value = "March"

# Row 34
$ws.Range("I34").Value = "January"

# Row 35
$ws.Range("A35").Value = "TO240_181_51_459 | Acropolis & 6 Archaeological Sites Combo Ticket"
$ws.Range("B35").Value = 18
$ws.Range("C35").Value = 0
$ws.Range("E35").Value = 18
$ws.Range("G35").Value = 36
$ws.Range("I35").Value = "February"

# Row 36
$ws.Range("A36").Value = "TO240_181_51_278 | Acropolis, Acropolis Museum & 6 Archaeological Sites Tickets"
$ws.Range("B36").Value = 12
$ws.Range("C36").Value = 6
$ws.Range("E36").Value = 0
$ws.Range("G36").Value = 18
$ws.Range("I36").Value = "March"

# Row 37
$ws.Range("A37").Value = "TO240_181_51_278 | Acropolis, Acropolis Museum & 6 Archaeological Sites Tickets"
$ws.Range("B37").Value = 12
$ws.Range("C37").Value = 6
$ws.Range("E37").Value = 0
$ws.Range("G37").Value = 18
$ws.Range("I37").Value = "June"

# Row 38
$ws.Range("I38").Value = "August"

# Row 39
$ws.Range("A39").Value = "TO240_181_51_459 | Acropolis & 6 Archaeological Sites Combo Ticket"
$ws.Range("B39").Value = 18
$ws.Range("C39").Value = 0
$ws.Range("E39").Value = 18
$ws.Range("G39").Value = 36
$ws.Range("I39").Value = "May"

# Row 45
$ws.Range("I45").Value = "June"
$ws.Range("K45").Value = 12

# Row 46
$ws.Range("I46").Value = "May"
$ws.Range("K46").Value = 20

# Row 47
$ws.Range("I47").Value = "July"
$ws.Range("K47").Value = 4
